$d = $word.ActiveDocument

# Paragraph 3 (0-index 2) is the "3. Sim, o SCRUM..." paragraph that also
# carries the trailing _GoBack bookmark at the very end of the story.
# We replace its whole Range with the expanded content: the original
# paragraph (now split into extra runs around "tem" with proofErr marks),
# four new paragraphs describing the team roles, and a final empty
# paragraph holding the _GoBack bookmark.

$p3 = $d.Paragraphs.Item(3)
$r = $p3.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00703010" w:rsidRPr="00342258" w:rsidRDefault="00703010"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">3. Sim, o SCRUM pode ser utilizado, por conta de se tratar de um método ágil, não precisando de todos os requisitos para começar o projeto, produz o projeto por partes acompanhando resultados para que seja possível diversas mudanças ao decorrer do projeto, sempre esperando por um feedback. Um método onde os funcionários escolhidos </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>tem</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> papéis definidos, mas são flexíveis, para uma realização mais rápida no Projeto.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>4. De acordo com o SCRUM, os profissionais terão papéis e responsabilidades detalhadas.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Um PO, existira vários produtos e soluções a entregar, e ele será responsável por isso.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Um Scrum Master, na qual pode ser dividido pelos times do Projeto.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Um time, o</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s profissionais do Time devem ter múltiplas e complementares competências para lidar com todas as tarefas de desenvolvimento</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$r.InsertXML($xml)
